$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 5).Value2 = 3.0
$ws.Cells.Item(2, 6).Value2 = 1.0
$ws.Cells.Item(2, 7).Value2 = 7.312497333333334
$ws.Cells.Item(2, 8).Value2 = 21.937492
$ws.Cells.Item(2, 9).Value2 = 0.05970572560549242
$ws.Cells.Item(2, 10).Value2 = 0.05970572560549242
$ws.Cells.Item(2, 11).Value2 = 3.0
$ws.Cells.Item(2, 12).Value2 = 1.0
$ws.Cells.Item(2, 13).Value2 = 98.946724
$ws.Cells.Item(2, 14).Value2 = 296.840172
$ws.Cells.Item(2, 15).Value2 = 0.2098009692989996
$ws.Cells.Item(2, 16).Value2 = 0.2098009692989996
$ws.Cells.Item(2, 17).Value2 = 723.5476553920694
$ws.Cells.Item(2, 18).Value2 = 6511.928898528625
$ws.Cells.Item(2, 19).Value2 = 0.01252631910473241
$ws.Cells.Item(2, 20).Value2 = 0.01252631910473241

$ws.Cells.Item(3, 5).Value2 = 3.0
$ws.Cells.Item(3, 6).Value2 = 1.0
$ws.Cells.Item(3, 7).Value2 = 7.312497333333334
$ws.Cells.Item(3, 8).Value2 = 21.937492
$ws.Cells.Item(3, 9).Value2 = 0.05970572560549242
$ws.Cells.Item(3, 10).Value2 = 0.05970572560549242
$ws.Cells.Item(3, 11).Value2 = 3.0
$ws.Cells.Item(3, 12).Value2 = 1.0
$ws.Cells.Item(3, 13).Value2 = 163.0062356666667
$ws.Cells.Item(3, 14).Value2 = 489.018707
$ws.Cells.Item(3, 15).Value2 = 0.345629090707923
$ws.Cells.Item(3, 16).Value2 = 0.3456290907079231
$ws.Cells.Item(3, 17).Value2 = 1191.982663629205
$ws.Cells.Item(3, 18).Value2 = 10727.84397266285
$ws.Cells.Item(3, 19).Value2 = 0.0206360356510831
$ws.Cells.Item(3, 20).Value2 = 0.0206360356510831

$ws.Cells.Item(4, 5).Value2 = 3.0
$ws.Cells.Item(4, 6).Value2 = 1.0
$ws.Cells.Item(4, 7).Value2 = 7.312497333333334
$ws.Cells.Item(4, 8).Value2 = 21.937492
$ws.Cells.Item(4, 9).Value2 = 0.05970572560549242
$ws.Cells.Item(4, 10).Value2 = 0.05970572560549242
$ws.Cells.Item(4, 11).Value2 = 3.0
$ws.Cells.Item(4, 12).Value2 = 1.0
$ws.Cells.Item(4, 13).Value2 = 65.39610666666668
$ws.Cells.Item(4, 14).Value2 = 196.18832
$ws.Cells.Item(4, 15).Value2 = 0.1386621609326595
$ws.Cells.Item(4, 16).Value2 = 0.1386621609326595
$ws.Cells.Item(4, 17).Value2 = 478.2088556103824
$ws.Cells.Item(4, 18).Value2 = 4303.879700493441
$ws.Cells.Item(4, 19).Value2 = 0.008278924932509996
$ws.Cells.Item(4, 20).Value2 = 0.008278924932509998

$ws.Cells.Item(5, 5).Value2 = 3.0
$ws.Cells.Item(5, 6).Value2 = 1.0
$ws.Cells.Item(5, 7).Value2 = 7.312497333333334
$ws.Cells.Item(5, 8).Value2 = 21.937492
$ws.Cells.Item(5, 9).Value2 = 0.05970572560549242
$ws.Cells.Item(5, 10).Value2 = 0.05970572560549242
$ws.Cells.Item(5, 11).Value2 = 3.0
$ws.Cells.Item(5, 12).Value2 = 1.0
$ws.Cells.Item(5, 13).Value2 = 144.2727966666667
$ws.Cells.Item(5, 14).Value2 = 432.81839
$ws.Cells.Item(5, 15).Value2 = 0.3059077790604178
$ws.Cells.Item(5, 16).Value2 = 0.3059077790604179
$ws.Cells.Item(5, 17).Value2 = 1054.994440897542
$ws.Cells.Item(5, 18).Value2 = 9494.94996807788
$ws.Cells.Item(5, 19).Value2 = 0.01826444591716691
$ws.Cells.Item(5, 20).Value2 = 0.01826444591716691

$ws.Cells.Item(6, 5).Value2 = 3.0
$ws.Cells.Item(6, 6).Value2 = 1.0
$ws.Cells.Item(6, 7).Value2 = 29.68221266666667
$ws.Cells.Item(6, 8).Value2 = 89.046638
$ws.Cells.Item(6, 9).Value2 = 0.242351957758873
$ws.Cells.Item(6, 10).Value2 = 0.242351957758873
$ws.Cells.Item(6, 11).Value2 = 3.0
$ws.Cells.Item(6, 12).Value2 = 1.0
$ws.Cells.Item(6, 13).Value2 = 98.946724
$ws.Cells.Item(6, 14).Value2 = 296.840172
$ws.Cells.Item(6, 15).Value2 = 0.2098009692989996
$ws.Cells.Item(6, 16).Value2 = 0.2098009692989996
$ws.Cells.Item(6, 17).Value2 = 2936.957704437971
$ws.Cells.Item(6, 18).Value2 = 26432.61933994174
$ws.Cells.Item(6, 19).Value2 = 0.05084567564932176
$ws.Cells.Item(6, 20).Value2 = 0.05084567564932176

$ws.Cells.Item(7, 5).Value2 = 3.0
$ws.Cells.Item(7, 6).Value2 = 1.0
$ws.Cells.Item(7, 7).Value2 = 29.68221266666667
$ws.Cells.Item(7, 8).Value2 = 89.046638
$ws.Cells.Item(7, 9).Value2 = 0.242351957758873
$ws.Cells.Item(7, 10).Value2 = 0.242351957758873
$ws.Cells.Item(7, 11).Value2 = 3.0
$ws.Cells.Item(7, 12).Value2 = 1.0
$ws.Cells.Item(7, 13).Value2 = 163.0062356666667
$ws.Cells.Item(7, 14).Value2 = 489.018707
$ws.Cells.Item(7, 15).Value2 = 0.345629090707923
$ws.Cells.Item(7, 16).Value2 = 0.3456290907079231
$ws.Cells.Item(7, 17).Value2 = 4838.385753050786
$ws.Cells.Item(7, 18).Value2 = 43545.47177745707
$ws.Cells.Item(7, 19).Value2 = 0.08376388679148425
$ws.Cells.Item(7, 20).Value2 = 0.08376388679148425

$ws.Cells.Item(8, 5).Value2 = 3.0
$ws.Cells.Item(8, 6).Value2 = 1.0
$ws.Cells.Item(8, 7).Value2 = 29.68221266666667
$ws.Cells.Item(8, 8).Value2 = 89.046638
$ws.Cells.Item(8, 9).Value2 = 0.242351957758873
$ws.Cells.Item(8, 10).Value2 = 0.242351957758873
$ws.Cells.Item(8, 11).Value2 = 3.0
$ws.Cells.Item(8, 12).Value2 = 1.0
$ws.Cells.Item(8, 13).Value2 = 65.39610666666668
$ws.Cells.Item(8, 14).Value2 = 196.18832
$ws.Cells.Item(8, 15).Value2 = 0.1386621609326595
$ws.Cells.Item(8, 16).Value2 = 0.1386621609326595
$ws.Cells.Item(8, 17).Value2 = 1941.101145652018
$ws.Cells.Item(8, 18).Value2 = 17469.91031086816
$ws.Cells.Item(8, 19).Value2 = 0.03360504616910594
$ws.Cells.Item(8, 20).Value2 = 0.03360504616910594

$ws.Cells.Item(9, 5).Value2 = 3.0
$ws.Cells.Item(9, 6).Value2 = 1.0
$ws.Cells.Item(9, 7).Value2 = 29.68221266666667
$ws.Cells.Item(9, 8).Value2 = 89.046638
$ws.Cells.Item(9, 9).Value2 = 0.242351957758873
$ws.Cells.Item(9, 10).Value2 = 0.242351957758873
$ws.Cells.Item(9, 11).Value2 = 3.0
$ws.Cells.Item(9, 12).Value2 = 1.0
$ws.Cells.Item(9, 13).Value2 = 144.2727966666667
$ws.Cells.Item(9, 14).Value2 = 432.81839
$ws.Cells.Item(9, 15).Value2 = 0.3059077790604178
$ws.Cells.Item(9, 16).Value2 = 0.3059077790604179
$ws.Cells.Item(9, 17).Value2 = 4282.335832674758
$ws.Cells.Item(9, 18).Value2 = 38541.02249407282
$ws.Cells.Item(9, 19).Value2 = 0.07413734914896103
$ws.Cells.Item(9, 20).Value2 = 0.07413734914896104

$ws.Cells.Item(10, 5).Value2 = 3.0
$ws.Cells.Item(10, 6).Value2 = 1.0
$ws.Cells.Item(10, 7).Value2 = 11.06470466666667
$ws.Cells.Item(10, 8).Value2 = 33.194114
$ws.Cells.Item(10, 9).Value2 = 0.09034208022509747
$ws.Cells.Item(10, 10).Value2 = 0.09034208022509749
$ws.Cells.Item(10, 11).Value2 = 3.0
$ws.Cells.Item(10, 12).Value2 = 1.0
$ws.Cells.Item(10, 13).Value2 = 98.946724
$ws.Cells.Item(10, 14).Value2 = 296.840172
$ws.Cells.Item(10, 15).Value2 = 0.2098009692989996
$ws.Cells.Item(10, 16).Value2 = 0.2098009692989996
$ws.Cells.Item(10, 17).Value2 = 1094.816278794179
$ws.Cells.Item(10, 18).Value2 = 9853.346509147608
$ws.Cells.Item(10, 19).Value2 = 0.01895385599971343
$ws.Cells.Item(10, 20).Value2 = 0.01895385599971344

$ws.Cells.Item(11, 5).Value2 = 3.0
$ws.Cells.Item(11, 6).Value2 = 1.0
$ws.Cells.Item(11, 7).Value2 = 11.06470466666667
$ws.Cells.Item(11, 8).Value2 = 33.194114
$ws.Cells.Item(11, 9).Value2 = 0.09034208022509747
$ws.Cells.Item(11, 10).Value2 = 0.09034208022509749
$ws.Cells.Item(11, 11).Value2 = 3.0
$ws.Cells.Item(11, 12).Value2 = 1.0
$ws.Cells.Item(11, 13).Value2 = 163.0062356666667
$ws.Cells.Item(11, 14).Value2 = 489.018707
$ws.Cells.Item(11, 15).Value2 = 0.345629090707923
$ws.Cells.Item(11, 16).Value2 = 0.3456290907079231
$ws.Cells.Item(11, 17).Value2 = 1803.615856476733
$ws.Cells.Item(11, 18).Value2 = 16232.5427082906
$ws.Cells.Item(11, 19).Value2 = 0.03122485104086267
$ws.Cells.Item(11, 20).Value2 = 0.03122485104086268

$ws.Cells.Item(12, 5).Value2 = 3.0
$ws.Cells.Item(12, 6).Value2 = 1.0
$ws.Cells.Item(12, 7).Value2 = 11.06470466666667
$ws.Cells.Item(12, 8).Value2 = 33.194114
$ws.Cells.Item(12, 9).Value2 = 0.09034208022509747
$ws.Cells.Item(12, 10).Value2 = 0.09034208022509749
$ws.Cells.Item(12, 11).Value2 = 3.0
$ws.Cells.Item(12, 12).Value2 = 1.0
$ws.Cells.Item(12, 13).Value2 = 65.39610666666668
$ws.Cells.Item(12, 14).Value2 = 196.18832
$ws.Cells.Item(12, 15).Value2 = 0.1386621609326595
$ws.Cells.Item(12, 16).Value2 = 0.1386621609326595
$ws.Cells.Item(12, 17).Value2 = 723.5886066164979
$ws.Cells.Item(12, 18).Value2 = 6512.297459548481
$ws.Cells.Item(12, 19).Value2 = 0.0125270280671637
$ws.Cells.Item(12, 20).Value2 = 0.0125270280671637

$ws.Cells.Item(13, 5).Value2 = 3.0
$ws.Cells.Item(13, 6).Value2 = 1.0
$ws.Cells.Item(13, 7).Value2 = 11.06470466666667
$ws.Cells.Item(13, 8).Value2 = 33.194114
$ws.Cells.Item(13, 9).Value2 = 0.09034208022509747
$ws.Cells.Item(13, 10).Value2 = 0.09034208022509749
$ws.Cells.Item(13, 11).Value2 = 3.0
$ws.Cells.Item(13, 12).Value2 = 1.0
$ws.Cells.Item(13, 13).Value2 = 144.2727966666667
$ws.Cells.Item(13, 14).Value2 = 432.81839
$ws.Cells.Item(13, 15).Value2 = 0.3059077790604178
$ws.Cells.Item(13, 16).Value2 = 0.3059077790604179
$ws.Cells.Item(13, 17).Value2 = 1596.335886550718
$ws.Cells.Item(13, 18).Value2 = 14367.02297895646
$ws.Cells.Item(13, 19).Value2 = 0.02763634511735766
$ws.Cells.Item(13, 20).Value2 = 0.02763634511735767

$ws.Cells.Item(14, 5).Value2 = 3.0
$ws.Cells.Item(14, 6).Value2 = 1.0
$ws.Cells.Item(14, 7).Value2 = 74.41623166666666
$ws.Cells.Item(14, 8).Value2 = 223.248695
$ws.Cells.Item(14, 9).Value2 = 0.6076002364105371
$ws.Cells.Item(14, 10).Value2 = 0.6076002364105371
$ws.Cells.Item(14, 11).Value2 = 3.0
$ws.Cells.Item(14, 12).Value2 = 1.0
$ws.Cells.Item(14, 13).Value2 = 98.946724
$ws.Cells.Item(14, 14).Value2 = 296.840172
$ws.Cells.Item(14, 15).Value2 = 0.2098009692989996
$ws.Cells.Item(14, 16).Value2 = 0.2098009692989996
$ws.Cells.Item(14, 17).Value2 = 7363.242335841726
$ws.Cells.Item(14, 18).Value2 = 66269.18102257553
$ws.Cells.Item(14, 19).Value2 = 0.127475118545232
$ws.Cells.Item(14, 20).Value2 = 0.127475118545232

$ws.Cells.Item(15, 5).Value2 = 3.0
$ws.Cells.Item(15, 6).Value2 = 1.0
$ws.Cells.Item(15, 7).Value2 = 74.41623166666666
$ws.Cells.Item(15, 8).Value2 = 223.248695
$ws.Cells.Item(15, 9).Value2 = 0.6076002364105371
$ws.Cells.Item(15, 10).Value2 = 0.6076002364105371
$ws.Cells.Item(15, 11).Value2 = 3.0
$ws.Cells.Item(15, 12).Value2 = 1.0
$ws.Cells.Item(15, 13).Value2 = 163.0062356666667
$ws.Cells.Item(15, 14).Value2 = 489.018707
$ws.Cells.Item(15, 15).Value2 = 0.345629090707923
$ws.Cells.Item(15, 16).Value2 = 0.3456290907079231
$ws.Cells.Item(15, 17).Value2 = 12130.30979648193
$ws.Cells.Item(15, 18).Value2 = 109172.7881683374
$ws.Cells.Item(15, 19).Value2 = 0.210004317224493
$ws.Cells.Item(15, 20).Value2 = 0.210004317224493

$ws.Cells.Item(16, 5).Value2 = 3.0
$ws.Cells.Item(16, 6).Value2 = 1.0
$ws.Cells.Item(16, 7).Value2 = 74.41623166666666
$ws.Cells.Item(16, 8).Value2 = 223.248695
$ws.Cells.Item(16, 9).Value2 = 0.6076002364105371
$ws.Cells.Item(16, 10).Value2 = 0.6076002364105371
$ws.Cells.Item(16, 11).Value2 = 3.0
$ws.Cells.Item(16, 12).Value2 = 1.0
$ws.Cells.Item(16, 13).Value2 = 65.39610666666668
$ws.Cells.Item(16, 14).Value2 = 196.18832
$ws.Cells.Item(16, 15).Value2 = 0.1386621609326595
$ws.Cells.Item(16, 16).Value2 = 0.1386621609326595
$ws.Cells.Item(16, 17).Value2 = 4866.531823804712
$ws.Cells.Item(16, 18).Value2 = 43798.78641424241
$ws.Cells.Item(16, 19).Value2 = 0.08425116176387983
$ws.Cells.Item(16, 20).Value2 = 0.08425116176387984

$ws.Cells.Item(17, 5).Value2 = 3.0
$ws.Cells.Item(17, 6).Value2 = 1.0
$ws.Cells.Item(17, 7).Value2 = 74.41623166666666
$ws.Cells.Item(17, 8).Value2 = 223.248695
$ws.Cells.Item(17, 9).Value2 = 0.6076002364105371
$ws.Cells.Item(17, 10).Value2 = 0.6076002364105371
$ws.Cells.Item(17, 11).Value2 = 3.0
$ws.Cells.Item(17, 12).Value2 = 1.0
$ws.Cells.Item(17, 13).Value2 = 144.2727966666667
$ws.Cells.Item(17, 14).Value2 = 432.81839
$ws.Cells.Item(17, 15).Value2 = 0.3059077790604178
$ws.Cells.Item(17, 16).Value2 = 0.3059077790604179
$ws.Cells.Item(17, 17).Value2 = 10736.23785994456
$ws.Cells.Item(17, 18).Value2 = 96626.14073950106
$ws.Cells.Item(17, 19).Value2 = 0.1858696388769322
$ws.Cells.Item(17, 20).Value2 = 0.1858696388769323

